# Work diary of the latest bugfix.
# Appends one new entry (row) to the "Journal de travail" table, describing
# a bugfix done on 2022-05-11 (serial date 44692):
#   Type        = Réalisation
#   Durée       = 2.5 heures
#   Description = BugFix du bug d'hier
#   Remarque    = prepare : mauvaise définiton des return dans la
#                 documentation de la fonction. Fonction non-définie : <?php

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel Table ("Tableau1") by one row - this also widens its
# `ref`/`autoFilter` range (A1:F43 -> A1:F44) automatically.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()
$r = $newListRow.Range.Row

# Fill in the new row's data.
$ws.Cells.Item($r, 1).Value = 44692
$ws.Cells.Item($r, 1).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item($r, 2).Value = "Réalisation"
$ws.Cells.Item($r, 3).Value = 2.5
$ws.Cells.Item($r, 4).Value = "BugFix du bug d'hier"
$ws.Cells.Item($r, 5).Value = "prepare : mauvaise définiton des return dans la documentation de la fonction. Fonction non-définie : <?php"

# Match the taller wrapped-text row height used by similar multi-line rows.
$ws.Rows.Item($r).RowHeight = 45
